$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06945541931069417
$ws.Range("H2").Value = -26.67359468159382
$ws.Range("I2").Value = 44.18294189445411
$ws.Range("G3").Value = 0.09349869261810012
$ws.Range("H3").Value = 41.46340837422618
$ws.Range("G4").Value = 0.0319872546868956
$ws.Range("H4").Value = 241.18764593717
$ws.Range("G5").Value = -0.002380750856928651
$ws.Range("H5").Value = 82.21968466619821
$ws.Range("G6").Value = -0.2224733353303241
$ws.Range("H6").Value = -0.5776095069149511
$ws.Range("G7").Value = -0.1632254436820946
$ws.Range("H7").Value = 34.67734261691146
$ws.Range("G8").Value = -0.4272931295765198
$ws.Range("H8").Value = -15.4048225799581
$ws.Range("G9").Value = -0.4677568661113735
$ws.Range("H9").Value = -17.33982358247239
$ws.Range("G10").Value = -0.04780475720886686
$ws.Range("H10").Value = -395.0020785111791
$ws.Range("G11").Value = 0.1006826925446078
$ws.Range("H11").Value = 726.6869411803648
$ws.Range("G12").Value = 0.2142891697627738
$ws.Range("H12").Value = -5.677934495648033
$ws.Range("G13").Value = 0.2575270914992639
$ws.Range("H13").Value = -2.207594958775476
$ws.Range("G14").Value = -0.04253632807769994
$ws.Range("H14").Value = -344.6024565808995
$ws.Range("G15").Value = 0.01217757038339743
$ws.Range("H15").Value = -39.67532997662359
$ws.Range("G16").Value = 0.1464452094648861
$ws.Range("H16").Value = 24.10705065002093
$ws.Range("G17").Value = 0.2072174633608206
$ws.Range("H17").Value = -5.307912052828935
$ws.Range("G18").Value = 0.03540522168747869
$ws.Range("H18").Value = -41.44665046469648
$ws.Range("G19").Value = 0.08896806781927279
$ws.Range("H19").Value = -1.24458313860303
$ws.Range("G20").Value = -0.1757920433900003
$ws.Range("H20").Value = -20.79355852941027
$ws.Range("G21").Value = -0.1495722484360507
$ws.Range("H21").Value = 25.15021371309935
$ws.Range("G22").Value = 0.03888876448762387
$ws.Range("H22").Value = -28.49784688123723
$ws.Range("G23").Value = 0.04882927976274308
$ws.Range("H23").Value = 19.56196417353189
$ws.Range("G24").Value = 0.1102964164904302
$ws.Range("H24").Value = -4.698611799212473
$ws.Range("G25").Value = 0.1564294645494417
$ws.Range("H25").Value = 2.862670742096184
$ws.Range("G26").Value = 0.005121160308161785
$ws.Range("H26").Value = -90.31491771793257
$ws.Range("G27").Value = 0.02571075827568861
$ws.Range("H27").Value = -49.05412179169874
$ws.Range("G28").Value = 0.1852169158322425
$ws.Range("H28").Value = 21.12720998834223
$ws.Range("G29").Value = 0.1712427901252587
$ws.Range("H29").Value = 0.3161642810312146
$ws.Range("G30").Value = 0.02012504302092714
$ws.Range("H30").Value = 2.856760009159318
$ws.Range("G31").Value = 0.03576073969376106
$ws.Range("H31").Value = 268.4785988493088
$ws.Range("G32").Value = 0.02306417761190661
$ws.Range("H32").Value = -38.15458877526477
$ws.Range("G33").Value = 0.007488741468394969
$ws.Range("H33").Value = -71.31150371297191
$ws.Range("G34").Value = 0.08028610370952603
$ws.Range("H34").Value = -37.26080090668251
$ws.Range("G35").Value = 0.1152670058733629
$ws.Range("H35").Value = -10.41001511111316
$ws.Range("G36").Value = -0.01544068208486425
$ws.Range("H36").Value = -202.7195162420344
$ws.Range("G37").Value = -0.01321741388607327
$ws.Range("H37").Value = -186.3065976376692
$ws.Range("G38").Value = -0.0002570172406311324
$ws.Range("H38").Value = 87.41600417728954
$ws.Range("G39").Value = -0.008343232672271031
$ws.Range("H39").Value = 75.0271494790854
$ws.Range("G40").Value = 0.1370862529706723
$ws.Range("H40").Value = -7.091027158794515
$ws.Range("G41").Value = 0.1454647524234101
$ws.Range("H41").Value = -9.872696445183125
$ws.Range("G42").Value = 0.03435261780493961
$ws.Range("H42").Value = -46.79366719130135
$ws.Range("G43").Value = 0.08121713682090591
$ws.Range("H43").Value = 133.6476248339604
$ws.Range("G44").Value = 0.03715818435271305
$ws.Range("H44").Value = 163.2945268371974
$ws.Range("G45").Value = 0.02236503426017778
$ws.Range("H45").Value = -45.52824468497757
$ws.Range("G46").Value = -0.06583910550071984
$ws.Range("H46").Value = -0.02923355798230415
$ws.Range("G47").Value = -0.08606578252273069
$ws.Range("H47").Value = -108.3427782513717
$ws.Range("G48").Value = -0.1033774318657585
$ws.Range("H48").Value = 17.93840374136519
$ws.Range("G49").Value = -0.1447366959655899
$ws.Range("H49").Value = 26.7086970526666
$ws.Range("G50").Value = 0.09686235988120949
$ws.Range("H50").Value = -11.03231878355261
$ws.Range("G51").Value = 0.1545542984444849
$ws.Range("H51").Value = 54.13685217849888
$ws.Range("G52").Value = 0.06859152454239748
$ws.Range("H52").Value = 15.04967083779549
$ws.Range("G53").Value = 0.05293799105957257
$ws.Range("H53").Value = -21.62290464147322
$ws.Range("G54").Value = -0.0922424933611274
$ws.Range("H54").Value = -31.92466785317035
$ws.Range("G55").Value = -0.0487036717490799
$ws.Range("H55").Value = 36.94191308394215
$ws.Range("G56").Value = 0.09211789837797493
$ws.Range("H56").Value = 101.0163453919548
$ws.Range("G57").Value = 0.1183627517465083
$ws.Range("H57").Value = 2189.291911497441
